# Fix pandas API compatibility issue in nz_util.py - recompute synthetic-control
# outcome/gap columns (C:E) on each sheet and drop the now-unused bold/bordered
# header/Year-column formatting.
$wb = $excel.ActiveWorkbook

# --- gdp_per_capita ---
$ws = $wb.Worksheets.Item("gdp_per_capita")
$colC = @(26556.58015712952,28347.25180397787,30011.87409579987,30785.90574021574,32762.96229035686,34678.72305017555,36912.16168001806,38172.02127435797,40985.66687073935,42064.27557717169,42181.37820617731,44075.73205593701,46006.20622525735,46070.52879091061,48437.44088493662,48784.65452206451,49412.78005855967,51952.00189555186,55755.39313506885,57735.96599455603,59396.35117674874,60905.88655061334,65343.15248542746)
$colD = @(47.41984287047671,-548.2518039778697,-511.8740957998707,109.0942597842622,-23.96229035685974,589.2769498244452,-226.1616800180636,339.9787256420313,12.33312926064536,-218.2755771716911,238.621793822691,547.2679440629872,1507.793774742648,3411.471209089388,5739.559115063377,6675.345477935494,6299.219941440329,4114.998104448139,3458.606864931149,3227.034005443973,3983.648823251264,1292.113449386663,2111.847514572539)
$colE = @(0.1785615564575852,-1.934056280901754,-1.705571915189085,0.3543643013294621,-0.07313835099676738,1.699246390854239,-0.6127023444971889,0.890648999691331,0.03009132265565326,-0.5189096309794752,0.5657041186666245,1.241653668663843,3.277370377727149,7.404888328007313,11.84942682809715,13.68328943462405,12.74815935062761,7.920769083588415,6.20317904772474,5.589295943793948,6.706891491359322,2.121491899330461,3.2319339276499)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $colC[$i]   # Synthetic
    $ws.Cells.Item($r, 4).Value = $colD[$i]   # Absolute Gap
    $ws.Cells.Item($r, 5).Value = $colE[$i]   # Percent Gap
}

# Drop the bold/bordered/centered style from the header row and Year column
# (matches upstream notebook re-export, which no longer applies that style)
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A2:A24").ClearFormats()

# --- population ---
$ws = $wb.Worksheets.Item("population")
$colC = @(484685.6265208999,489816.3383665308,498545.7708089837,511948.5609111016,524050.6765385518,532793.9584914244,540777.7451329936,548194.7219745264,553870.2957903675,559156.9886639647,565975.1799374919,572932.6545736709,578995.8276273068,584190.1940168024,592115.0299842544,603944.4211149375,617287.0938656105,631410.3692954445,644192.5032100378,655322.6097235683,668009.5945787356,673787.3964828703,671726.9786689661)
$colD = @(8814.373479100119,6383.661633469164,3154.229191016289,-748.5609111016383,-2850.67653855175,-3393.958491424448,-3877.745132993557,-3394.721974526416,-2370.295790367527,-1356.988663964672,-675.1799374918919,-8332.654573670938,-20495.8276273068,-24090.19401680236,-24815.02998425439,-25044.42111493752,-22587.09386561054,-22210.36929544446,-24592.50321003783,-27122.60972356831,-28609.59457873565,-26387.39648287033,-19726.9786689661)
$colE = @(1.818575380988741,1.303276582148686,0.6326859790421977,-0.1462179930283315,-0.5439696323609345,-0.6370114445430738,-0.7170681796529742,-0.6192547717896785,-0.4279514190204293,-0.2426847363934457,-0.1192949728937691,-1.454386393785044,-3.539892111364185,-4.123690240529693,-4.190913712309292,-4.14680891806289,-3.659090573911784,-3.517580701157595,-3.817570537920322,-4.138817938085383,-4.282811925295416,-3.91627932202516,-2.936755452052158)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $colC[$i]   # Synthetic
    $ws.Cells.Item($r, 4).Value = $colD[$i]   # Absolute Gap
    $ws.Cells.Item($r, 5).Value = $colE[$i]   # Percent Gap
}

# Drop the bold/bordered/centered style from the header row and Year column
# (matches upstream notebook re-export, which no longer applies that style)
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A2:A24").ClearFormats()

# --- total_gdp ---
$ws = $wb.Worksheets.Item("total_gdp")
$colC = @(13385.91504048598,13885.62249731768,14943.84358210379,16165.41077190083,17421.60624530239,18785.61193577181,19720.63996221944,20535.98000526109,22058.13548824816,21903.10623419921,22622.32181288475,23664.38036372572,24867.33212940194,25634.10556249695,27101.50915744621,28924.90985626361,31367.33230213309,33515.44365851385,35750.83377560642,37747.53732496355,39609.98242851395,40481.51672776663,44395.46659297859)
$colD = @(-256.9150404859811,-91.62249731768316,-143.8435821037929,-371.4107719008316,-357.6062453023878,-114.6119357718089,-23.63996221944399,445.0199947389119,552.864511751839,1438.893765800793,1357.678187115253,1529.619636274278,1669.667870598063,2080.894437503051,3633.490842553794,3181.090143736386,1764.667697866913,640.5563414861463,938.1662243935789,549.4626750364478,915.0175714860525,-214.5167277666333,-414.4665929785915)
$colE = @(-1.919293822715415,-0.6598371613183498,-0.9625608118386272,-2.297564702447451,-2.052659440623127,-0.6101048832674084,-0.1198742143497023,2.167025847438996,2.506397297479594,6.569359388642894,6.001497982147768,6.463805993496353,6.714302370313091,8.117679130366962,13.40696867264856,10.99775300785433,5.625813763406681,1.911227397174637,2.624179984953834,1.455625224782736,2.310068107546952,-0.5299127728074832,-0.9335786394103174)
for ($i = 0; $i -lt $colC.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $colC[$i]   # Synthetic
    $ws.Cells.Item($r, 4).Value = $colD[$i]   # Absolute Gap
    $ws.Cells.Item($r, 5).Value = $colE[$i]   # Percent Gap
}

# Drop the bold/bordered/centered style from the header row and Year column
# (matches upstream notebook re-export, which no longer applies that style)
$ws.Range("A1:E1").ClearFormats()
$ws.Range("A2:A24").ClearFormats()

